$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add a new "2021" data column (column R) for rows 4-44, mirroring the
#    formatting already used by column Q (the 2020 column) for each row.
# ---------------------------------------------------------------------------

# Copy the per-row formatting from column Q into column R in one bulk
# operation so every row in R4:R44 picks up the same style Q uses on that row.
$ws.Range("Q4:Q44").Copy() | Out-Null
$ws.Range("R4:R44").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Fill in the values for the new column (2021 figures).
# ---------------------------------------------------------------------------

$ws.Range("R4").Value = 2021

$ws.Range("R5").Value = 33.262233298138462
# R6 stays blank (header/spacer row)

$ws.Range("R7").Value = 33.10733359588
$ws.Range("R8").Value = 33.434791912906064
# R9 stays blank (spacer row)
$ws.Range("R10").Value = 33.257744153601877
$ws.Range("R11").Value = 33.264901349775037
# R12 stays blank (spacer row)
$ws.Range("R13").Value = 40.506761643955343
$ws.Range("R14").Value = 31.172874511993772
$ws.Range("R15").Value = 29.650854888195603
$ws.Range("R16").Value = 24.097659040517694
# R17 stays blank (spacer row)
$ws.Range("R18").Value = 40.671537678982844
$ws.Range("R19").Value = 40.835744431591088
$ws.Range("R20").Value = 40.494322790314847
$ws.Range("R21").Value = 43.213901910043809
$ws.Range("R22").Value = 43.25339669708363
$ws.Range("R23").Value = 43.172054336673064
$ws.Range("R24").Value = 38.130315382405762
$ws.Range("R25").Value = 39.052986923894757
$ws.Range("R26").Value = 37.162402419999459
$ws.Range("R27").Value = 39.151777291250745
$ws.Range("R28").Value = 39.257524172776719
$ws.Range("R29").Value = 39.034895273433577
$ws.Range("R30").Value = 23.758244663001044
$ws.Range("R31").Value = 23.600989554960133
$ws.Range("R32").Value = 23.928051635532135
$ws.Range("R33").Value = 23.479873624436866
$ws.Range("R34").Value = 24.191255182099319
$ws.Range("R35").Value = 22.759796314587014
$ws.Range("R36").Value = 27.00352039684709
$ws.Range("R37").Value = 26.49991066711625
$ws.Range("R38").Value = 27.572066513890491
$ws.Range("R39").Value = 35.769118603355516
$ws.Range("R40").Value = 35.037028356467729
$ws.Range("R41").Value = 36.714712776303102
$ws.Range("R42").Value = 28.585005142864613
$ws.Range("R43").Value = 27.745516596290607
$ws.Range("R44").Value = 29.549901683892372

# ---------------------------------------------------------------------------
# 3. Update the active selection shown in the sheet view.
# ---------------------------------------------------------------------------
[void]$ws.Range("O10").Select()
